$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.164.88'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '2.323.17'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''303.17'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("D6").Value = '''99.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("D7").Value = '''0.508'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +2.31%  '
$ws.Range("D10").Value = '''36.12'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.76%  '
$ws.Range("D11").Value = '''0.0794'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("D13").Value = '''17.79'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.49%  '
$ws.Range("D14").Value = '''6.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.34%  '
$ws.Range("D15").Value = '2.683.65'
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("D16").Value = '2.301.96'
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("D17").Value = '''0.799'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.03%  '
$ws.Range("D18").Value = '43.094.11'
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("D19").Value = '''12.94'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.50%  '
$ws.Range("D20").Value = '''6.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.42%  '
$ws.Range("E21").Value = '  +1.05%  '
$ws.Range("D22").Value = '''68.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.71%  '
$ws.Range("D23").Value = '''240.14'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("E24").Value = '  -0.75%  '
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '''25.54'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.21%  '
$ws.Range("D28").Value = '''168.15'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").Value = '''34.41'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.93%  '
$ws.Range("D30").Value = '''9.20'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.77%  '
$ws.Range("E31").Value = '  -10.63%  '
$ws.Range("B32").Value = 'RenderToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D32").Value = '''4.97'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.01%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''5.17'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.47%  '
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").Value = '''0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").Value = '''17.79'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.17%  '
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("D37").Value = '''0.0699'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.07%  '
$ws.Range("E38").Value = '  +2.23%  '
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("E40").Value = '  -0.34%  '
$ws.Range("E41").Value = '  +0.30%  '
$ws.Range("D42").Value = '1.995.53'
$ws.Range("E42").Value = '  +0.23%  '
$ws.Range("E43").Value = '  +1.69%  '
$ws.Range("D44").Value = '''2.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.03%  '
$ws.Range("D45").Value = '''10.11'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.09%  '
$ws.Range("D46").Value = '''17.66'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").Value = '''76.55'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.48%  '
$ws.Range("D49").Value = '''55.21'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.31%  '
$ws.Range("D50").Value = '2.548.79'
$ws.Range("E50").Value = '  +0.84%  '
$ws.Range("D51").Value = '''1.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.48%  '
